$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 720; this shifts the existing rows
# 720-841 down to 722-843 (preserving their content/formatting) and
# leaves two blank rows at 720-721 for the new weekly records.
$ws.Rows.Item(720).Insert()
$ws.Rows.Item(720).Insert()

# New record #1 (row 720) - same dimensions as the old row 720 record,
# new date / volume / price.
$ws.Cells.Item(720, 1).Value = 9
$ws.Cells.Item(720, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(720, 3).Value = "Metropolitana"
$ws.Cells.Item(720, 4).Value = 44951
$ws.Cells.Item(720, 5).Value = 13
$ws.Cells.Item(720, 6).Value = 100112040
$ws.Cells.Item(720, 7).Value = "Cilantro"
$ws.Cells.Item(720, 8).Value = "Sin especificar"
$ws.Cells.Item(720, 9).Value = "Primera"
$ws.Cells.Item(720, 10).Value = 70
$ws.Cells.Item(720, 11).Value = 12000
$ws.Cells.Item(720, 12).Value = 12000
$ws.Cells.Item(720, 13).Value = 12000
$ws.Cells.Item(720, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(720, 15).Value = "Región Metropolitana"
$ws.Cells.Item(720, 16).Value = 333
$ws.Cells.Item(720, 17).Value = 36
$ws.Cells.Item(720, 18).Value = "Hortaliza"

# New record #2 (row 721)
$ws.Cells.Item(721, 1).Value = 9
$ws.Cells.Item(721, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(721, 3).Value = "Metropolitana"
$ws.Cells.Item(721, 4).Value = 44951
$ws.Cells.Item(721, 5).Value = 13
$ws.Cells.Item(721, 6).Value = 100112040
$ws.Cells.Item(721, 7).Value = "Cilantro"
$ws.Cells.Item(721, 8).Value = "Sin especificar"
$ws.Cells.Item(721, 9).Value = "Primera"
$ws.Cells.Item(721, 10).Value = 160
$ws.Cells.Item(721, 11).Value = 24000
$ws.Cells.Item(721, 12).Value = 25000
$ws.Cells.Item(721, 13).Value = 24500
$ws.Cells.Item(721, 14).Value = "$/docena de atados"
$ws.Cells.Item(721, 15).Value = "Región Metropolitana"
$ws.Cells.Item(721, 16).Value = 8167
$ws.Cells.Item(721, 17).Value = 3
$ws.Cells.Item(721, 18).Value = "Hortaliza"
